# Insert a new data row at row 254 (pushing existing rows 254:351 down to 255:352)
# and populate it with the new Ají record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(254).Insert()

$ws.Cells.Item(254, 1).Value  = 5
$ws.Cells.Item(254, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(254, 3).Value  = "Maule"
$ws.Cells.Item(254, 4).Value  = 45009
$ws.Cells.Item(254, 5).Value  = 7
$ws.Cells.Item(254, 6).Value  = 100112021
$ws.Cells.Item(254, 7).Value  = "Ají"
$ws.Cells.Item(254, 8).Value  = "Cristal"
$ws.Cells.Item(254, 9).Value  = "Primera"
$ws.Cells.Item(254, 10).Value = 150
$ws.Cells.Item(254, 11).Value = 12000
$ws.Cells.Item(254, 12).Value = 12000
$ws.Cells.Item(254, 13).Value = 12000
$ws.Cells.Item(254, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(254, 15).Value = "Región del Maule"
$ws.Cells.Item(254, 16).Value = 480
$ws.Cells.Item(254, 17).Value = 25
$ws.Cells.Item(254, 18).Value = "Hortaliza"
